$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the Handedness (column C) for the five pitchers that were
# already present in rows 172-176 but missing their handedness value ---
$ws.Range("C172").Value = "LHP"
$ws.Range("C173").Value = "LHP"
$ws.Range("C174").Value = "RHP"
$ws.Range("C175").Value = "RHP"
$ws.Range("C176").Value = "RHP"

# --- Append the five new pitchers (rows 177-181) ---
# Row 177: Dane Dunning
$ws.Range("A177").Value = "Dane Dunning"
$ws.Range("B177").Value = "https://midfield.mlbstatic.com/v1/people/641540/spots/120"

# Row 178: Andre Pallante
$ws.Range("A178").Value = "Andre Pallante"

# Row 179: Keider Montero
$ws.Range("A179").Value = "Keider Montero"

# Row 180: Spencer Schwellenbach
$ws.Range("A180").Value = "Spencer Schwellenbach"

$ws.Range("B178").Value = "https://midfield.mlbstatic.com/v1/people/669467/spots/120"
$ws.Range("B179").Value = "https://midfield.mlbstatic.com/v1/people/672456/spots/120"
$ws.Range("B180").Value = "https://midfield.mlbstatic.com/v1/people/680885/spots/120"

# Row 181: David Peterson (no handedness recorded yet)
$ws.Range("A181").Value = "David Peterson"
$ws.Range("B181").Value = "https://midfield.mlbstatic.com/v1/people/656849/spots/120"

# Handedness for the new rows with known values
$ws.Range("C177").Value = "RHP"
$ws.Range("C178").Value = "RHP"
$ws.Range("C179").Value = "RHP"
$ws.Range("C180").Value = "RHP"

# --- Wire up the headshot hyperlinks for the new rows & restyle them ---
$ws.Hyperlinks.Add($ws.Range("B177"), "https://midfield.mlbstatic.com/v1/people/641540/spots/120")
$ws.Hyperlinks.Add($ws.Range("B178"), "https://midfield.mlbstatic.com/v1/people/669467/spots/120")
$ws.Hyperlinks.Add($ws.Range("B179"), "https://midfield.mlbstatic.com/v1/people/672456/spots/120")
$ws.Hyperlinks.Add($ws.Range("B180"), "https://midfield.mlbstatic.com/v1/people/680885/spots/120")
$ws.Hyperlinks.Add($ws.Range("B181"), "https://midfield.mlbstatic.com/v1/people/656849/spots/120")

$ws.Range("B177").Style = "Hyperlink"
$ws.Range("B178").Style = "Hyperlink"
$ws.Range("B179").Style = "Hyperlink"
$ws.Range("B180").Style = "Hyperlink"
$ws.Range("B181").Style = "Hyperlink"

# --- Column B is no longer the widest-possible fixed width; narrow it to
# (approximately) fit its new best-fit content width ---
$ws.Columns.Item(2).ColumnWidth = 49.6

# --- Restore the view state (scrolled/selected cell moved down as rows were added) ---
$ws.Range("B182").Select()
